$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 181, shifting rows 181:256 down to 182:257
$ws.Rows.Item(181).Insert()

# Populate the new row 181 with the new data
$ws.Cells.Item(181, 1).Value = 3
$ws.Cells.Item(181, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(181, 3).Value = "Coquimbo"
$ws.Cells.Item(181, 4).Value = 44636
$ws.Cells.Item(181, 4).NumberFormat = $ws.Cells.Item(182, 4).NumberFormat
$ws.Cells.Item(181, 5).Value = 5
$ws.Cells.Item(181, 6).Value = 100112001
$ws.Cells.Item(181, 7).Value = "Berenjena"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 120
$ws.Cells.Item(181, 11).Value = 8500
$ws.Cells.Item(181, 12).Value = 9000
$ws.Cells.Item(181, 13).Value = 8750
$ws.Cells.Item(181, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(181, 15).Value = "Región Metropolitana"
$ws.Cells.Item(181, 16).Value = 146
$ws.Cells.Item(181, 17).Value = 60
$ws.Cells.Item(181, 18).Value = "Hortaliza"
